$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 411, pushing existing rows 411-480 down to 412-481.
$ws.Rows.Item(411).Insert()

# Populate the newly inserted row 411 with the new weekly data entry.
$ws.Cells.Item(411, 1).Value = 4
$ws.Cells.Item(411, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(411, 3).Value = "Los Lagos"
$ws.Cells.Item(411, 4).Value = 45218
$ws.Cells.Item(411, 4).NumberFormat = $ws.Cells.Item(412, 4).NumberFormat
$ws.Cells.Item(411, 5).Value = 10
$ws.Cells.Item(411, 6).Value = 100112032
$ws.Cells.Item(411, 7).Value = "Zapallo italiano"
$ws.Cells.Item(411, 8).Value = "Sin especificar"
$ws.Cells.Item(411, 9).Value = "Primera"
$ws.Cells.Item(411, 10).Value = 140
$ws.Cells.Item(411, 11).Value = 22000
$ws.Cells.Item(411, 12).Value = 22000
$ws.Cells.Item(411, 13).Value = 22000
$ws.Cells.Item(411, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(411, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(411, 16).Value = 440
$ws.Cells.Item(411, 17).Value = 50
$ws.Cells.Item(411, 18).Value = "Hortaliza"
